$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting rows 28-30 down to 29-31.
# This also auto-expands the SUM(...) ranges in D15:H15 and K30/K31's
# total formula because row 28 falls inside the summed 25:29 ranges.
$ws.Rows("28:28").Insert()

# Row 27: task label + Thurs value changed (3 -> 1), now referencing a
# brand-new shared string "24-Dann" (previously "24a-Dann").
$ws.Range("J27").Value = "24-Dann"
$ws.Range("L27").Value = 1

# New row 28: "24a-Dann" task with Sunday value 3.
$ws.Range("J28").Value = "24a-Dann"
$ws.Range("O28").Value = 3

# Row 29 (previously row 28, shifted down by the insert): keep the
# "24b-Dann" label, but move its value from Thurs (L) to Sun (O) = 3.
$ws.Range("J29").Value = "24b-Dann"
$ws.Range("L29").Clear()
$ws.Range("O29").Value = 3

# The row-insert also copied the row-above's formatting into a couple of
# now-unused cells (Thurs column on the new row, and the trailing blank
# "Sun-note" placeholder) - drop those so the row matches the source data.
$ws.Range("L28").Clear()
$ws.Range("Q28").Clear()

# Update the sheet view to match the new active cell/scroll position.
$ws.Range("I13").Select()
